$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 7 ("System Usability") - Content Placeholder 2
#   - "... The system should have:" -> split the trailing words and change
#     "should" to "shall": "... The " + "system shall have:"
#   - "Graphic User Interface (GUI). " -> split into "Graphic " + "User
#     Interface (GUI). "
#   - "Web enabled front end." -> split into "Web " + "enabled front end."
# ---------------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$body7 = $slide7.Shapes.Item(2).TextFrame.TextRange

$para1 = $body7.Paragraphs(1, 1)
$para1Text = $para1.Text
$marker1 = "system should have:"
$pos1 = $para1Text.IndexOf($marker1)
$run1b = $para1.Characters($pos1 + 1, $marker1.Length)
$run1b.Text = "system shall have:"

$body7b = $slide7.Shapes.Item(2).TextFrame.TextRange
$para2 = $body7b.Paragraphs(2, 1)
$split2 = "Graphic "
$run2a = $para2.Characters(1, $split2.Length)
$run2a.Text = $split2

$body7c = $slide7.Shapes.Item(2).TextFrame.TextRange
$para3 = $body7c.Paragraphs(3, 1)
$split3 = "Web "
$run3a = $para3.Characters(1, $split3.Length)
$run3a.Text = $split3

# ---------------------------------------------------------------------------
# Slide 8 ("User Roles and Accessibility") - Content Placeholder 2
#   The placeholder is empty (just an endParaRPr). Fill it in with the three
#   user-role paragraphs.
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)
$body8 = $slide8.Shapes.Item(2).TextFrame.TextRange

$dash = [char]0x2013

$localUserText = "Local User " + $dash + " user with locally installed test engine on PC/laptop.  Able to register with Test Server to register test engine in database.  Able to run tests locally.  Uses UI with internal engine."
$remoteUserText = "Remote User " + $dash + " user with same capability as local user, but has ability to view available test engines across the web and to use other hardware, servers,  infrastructure for test purposes as well as capability to view archived results."
$adminText = "Administrator " + $dash + " user with the ability to install application remotely, update application services, add remote users, de-register test engines, perform HA/DR testing.  "

$runLocal = $body8.InsertAfter($localUserText)
$runLocal.LanguageID = "en-US"
$runLocal.InsertAfter("`r")

$body8b = $slide8.Shapes.Item(2).TextFrame.TextRange
$runRemote = $body8b.InsertAfter($remoteUserText)
$runRemote.LanguageID = "en-US"
$runRemote.InsertAfter("`r")

$body8c = $slide8.Shapes.Item(2).TextFrame.TextRange
$runAdmin = $body8c.InsertAfter($adminText)
$runAdmin.LanguageID = "en-US"

# Split out the single "t" of "Able to run" in paragraph 1 into its own run,
# matching "Able " + "t" + "o run tests locally...".
$body8d = $slide8.Shapes.Item(2).TextFrame.TextRange
$firstPara = $body8d.Paragraphs(1, 1)
$firstParaText = $firstPara.Text
$ableMarker = "Able "
$splitPos = $firstParaText.LastIndexOf($ableMarker) + $ableMarker.Length
$tRun = $firstPara.Characters($splitPos + 1, 1)
$tRun.Text = "t"
$tRun.LanguageID = "en-US"
